# NESP-GET converted to RDF/TTL
# Replace the literal "bath:..." labels in column A of the SSSOM sheet with
# formulas derived from column B, and populate column O with a descriptive
# "mapping to IUCN GET" label formula, for rows 2-5.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("SSSOM")

$ws.Range("A2").Formula = '=_xlfn.CONCAT("bath:",LEFT(B2,FIND(" ",B2)-1))'
$ws.Range("A3").Formula = '=_xlfn.CONCAT("bath:",LEFT(B3,FIND(" ",B3)-1))'
$ws.Range("A4").Formula = '=_xlfn.CONCAT("bath:",LEFT(B4,FIND(" ",B4)-1))'
$ws.Range("A5").Formula = '=_xlfn.CONCAT("bath:",LEFT(B5,FIND(" ",B5)-1))'

$ws.Range("O2").Formula = '=_xlfn.CONCAT(B2, " - mapping to IUCN GET - ", ROW(B2)-1)'
$ws.Range("O3").Formula = '=_xlfn.CONCAT(B3, " - mapping to IUCN GET - ", ROW(B3)-1)'
$ws.Range("O4").Formula = '=_xlfn.CONCAT(B4, " - mapping to IUCN GET - ", ROW(B4)-1)'
$ws.Range("O5").Formula = '=_xlfn.CONCAT(B5, " - mapping to IUCN GET - ", ROW(B5)-1)'

# Make SSSOM the active/selected sheet with O3:O5 selected (matches the
# author's final cursor position when they finished curating the sheet).
$ws.Activate() | Out-Null
$ws.Range("O3:O5").Select() | Out-Null
